# Updates crypto price/volume data per upstream diff (Fri Aug  9 15:37:25 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.195.58'
$ws.Range("E2").Value = '  +2.82%  '
$ws.Range("D3").Value = '2.566.85'
$ws.Range("E3").Value = '  +4.49%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.42'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.95%  '
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("E8").Value = '  -6.06%  '
$ws.Range("D9").Value = '2.583.03'
$ws.Range("E9").Value = '  +3.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.28%  '
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("E12").Value = '  +1.46%  '
$ws.Range("E13").Value = '  +0.58%  '
$ws.Range("D14").Value = '3.023.51'
$ws.Range("E14").Value = '  +4.89%  '
$ws.Range("D15").Value = '60.236.30'
$ws.Range("E15").Value = '  +3.10%  '
$ws.Range("E16").Value = '  +0.34%  '
$ws.Range("E17").Value = '  +2.52%  '
$ws.Range("D18").Value = '2.581.72'
$ws.Range("E18").Value = '  +4.01%  '
$ws.Range("E19").Value = '  +1.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '346.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.43%  '
$ws.Range("E21").Value = '  +1.40%  '
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.29%  '
$ws.Range("E25").Value = '  +2.11%  '
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("D27").Value = '2.687.95'
$ws.Range("E27").Value = '  +4.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.993'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.23%  '
$ws.Range("D29").Value = '0.0₃0847'
$ws.Range("E29").Value = '  +4.98%  '
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.14'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.56%  '
$ws.Range("E34").Value = '  +1.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.69'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.61%  '
$ws.Range("E36").Value = '  +3.21%  '
$ws.Range("E37").Value = '  +1.26%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.859'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +21.64%  '
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.843'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.76'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.21%  '
$ws.Range("E41").Value = '  +2.80%  '
$ws.Range("E42").Value = '  +4.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '35.49'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0562'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.82%  '
$ws.Range("E45").Value = '  +1.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0994'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.75%  '
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.84%  '
$ws.Range("E50").Value = '  -0.94%  '
$ws.Range("D51").Value = '2.017.61'
$ws.Range("E51").Value = '  +5.56%  '
